$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$headerRange = $ws.Range("A1:U1")
$headerRange.Style = "Normal"
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U57"), [System.Reflection.Missing]::Value, 1)
$lo.TableStyle.Name = ""
